$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test app name and description values
$ws.Range("C2").Value = "TestApp1"
$ws.Range("D2").Value = "Test the creation of an app."

# Update the active cell selection
$ws.Range("A5").Select()
